$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename: split "name" into a clearer "team" / "team member" naming scheme ---
# A1: "team"  -> "team_type"
# B1: "name"  -> "team_member_name"
$ws.Range("A1").Value = "team_type"
$ws.Range("B1").Value = "team_member_name"

# Column B needs to be noticeably wider now that the header text is longer.
$ws.Columns("B").ColumnWidth = 22

# Row heights were re-flowed by Excel after the edit (wrap-text autofit).
$ws.Rows(3).RowHeight = 202.8
$ws.Rows(4).RowHeight = 72
$ws.Rows(5).RowHeight = 78
$ws.Rows(7).RowHeight = 156
$ws.Rows(8).RowHeight = 171.6
$ws.Rows(9).RowHeight = 202.8
$ws.Rows(10).RowHeight = 124.8
$ws.Rows(11).RowHeight = 374.4
$ws.Rows(12).RowHeight = 280.8
$ws.Rows(13).RowHeight = 124.8
$ws.Rows(14).RowHeight = 218.4
$ws.Rows(15).RowHeight = 93.6

# Update the window scroll/selection to reflect where the user was working.
$ws.Range("D14").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
